$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "B2" was the merged header's sub-label for the "total" column; the
# corrected data dropped the stray "unnamed: 1_level_1" label and now
# just repeats "total" there.
$ws.Range("B2").Value = "total"

# Two label-only rows ("situação do domicílio" and "grandes regiões e
# unidades da federação") were removed; deleting the lower one first
# keeps the row numbers of the other deletion valid.
$ws.Rows(8).Delete()
$ws.Rows(5).Delete()
